$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 30,5

$arr[0,0] = 0.6290000081062317
$arr[0,1] = 0.627177700348432
$arr[0,2] = 0.5326032852165257
$arr[0,3] = 0.6291687406669985
$arr[0,4] = 0.5286212045793928
$arr[1,0] = 0.6445000171661377
$arr[1,1] = 0.6450970632155301
$arr[1,2] = 0.5340965654554505
$arr[1,3] = 0.6445993031358885
$arr[1,4] = 0.5425584868093579
$arr[2,0] = 0.6380000114440918
$arr[2,1] = 0.6386261821801892
$arr[2,2] = 0.5460428073668492
$arr[2,3] = 0.640119462419114
$arr[2,4] = 0.5375808860129417
$arr[3,0] = 0.6345000267028809
$arr[3,1] = 0.6331508213041314
$arr[3,2] = 0.5316077650572424
$arr[3,3] = 0.6346441015430563
$arr[3,4] = 0.5151816824290691
$arr[4,0] = 0.6355000138282776
$arr[4,1] = 0.6376306620209059
$arr[4,2] = 0.5306122448979592
$arr[4,3] = 0.6346441015430563
$arr[4,4] = 0.5226480836236934
$arr[5,0] = 0.6430000066757202
$arr[5,1] = 0.6441015430562469
$arr[5,2] = 0.5281234444997511
$arr[5,3] = 0.6436037829766053
$arr[5,4] = 0.5211548033847685
$arr[6,0] = 0.6395000219345093
$arr[6,1] = 0.6406172224987556
$arr[6,2] = 0.5530114484818317
$arr[6,3] = 0.6411149825783972
$arr[6,4] = 0.560477849676456
$arr[7,0] = 0.6430000066757202
$arr[7,1] = 0.6455948232951717
$arr[7,2] = 0.5316077650572424
$arr[7,3] = 0.6441015430562469
$arr[7,4] = 0.5435540069686411
$arr[8,0] = 0.640999972820282
$arr[8,1] = 0.6411149825783972
$arr[8,2] = 0.5151816824290691
$arr[8,3] = 0.6416127426580388
$arr[8,4] = 0.5360876057740169
$arr[9,0] = 0.6499999761581421
$arr[9,1] = 0.6485813837730214
$arr[9,2] = 0.5231458437033349
$arr[9,3] = 0.6495769039323046
$arr[9,4] = 0.547536087605774
$arr[10,0] = 0.6359999775886536
$arr[10,1] = 0.6376306620209059
$arr[10,2] = 0.5266301642608263
$arr[10,3] = 0.6366351418616227
$arr[10,4] = 0.5186660029865605
$arr[11,0] = 0.6420000195503235
$arr[11,1] = 0.6416127426580388
$arr[11,2] = 0.5072175211548033
$arr[11,3] = 0.6421105027376804
$arr[11,4] = 0.538576406172225
$arr[12,0] = 0.6414999961853027
$arr[12,1] = 0.6411149825783972
$arr[12,2] = 0.5311100049776007
$arr[12,3] = 0.6411149825783972
$arr[12,4] = 0.5156794425087108
$arr[13,0] = 0.6579999923706055
$arr[13,1] = 0.658038825286212
$arr[13,2] = 0.5176704828272772
$arr[13,3] = 0.6570433051269288
$arr[13,4] = 0.5440517670482827
$arr[14,0] = 0.6365000009536743
$arr[14,1] = 0.6366351418616227
$arr[14,2] = 0.5629666500746641
$arr[14,3] = 0.6361373817819811
$arr[14,4] = 0.5559980089596814
$arr[15,0] = 0.6334999799728394
$arr[15,1] = 0.6326530612244898
$arr[15,2] = 0.5246391239422599
$arr[15,3] = 0.6331508213041314
$arr[15,4] = 0.5311100049776007
$arr[16,0] = 0.6359999775886536
$arr[16,1] = 0.6381284221005475
$arr[16,2] = 0.5316077650572424
$arr[16,3] = 0.6381284221005475
$arr[16,4] = 0.5246391239422599
$arr[17,0] = 0.6414999961853027
$arr[17,1] = 0.6406172224987556
$arr[17,2] = 0.538576406172225
$arr[17,3] = 0.6416127426580388
$arr[17,4] = 0.5156794425087108
$arr[18,0] = 0.6424999833106995
$arr[18,1] = 0.6426082628173221
$arr[18,2] = 0.5261324041811847
$arr[18,3] = 0.6426082628173221
$arr[18,4] = 0.5271279243404678
$arr[19,0] = 0.6420000195503235
$arr[19,1] = 0.6436037829766053
$arr[19,2] = 0.5350920856147336
$arr[19,3] = 0.6426082628173221
$arr[19,4] = 0.5271279243404678
$arr[20,0] = 0.6324999928474426
$arr[20,1] = 0.6306620209059234
$arr[20,2] = 0.5286212045793928
$arr[20,3] = 0.6321553011448482
$arr[20,4] = 0.5311100049776007
$arr[21,0] = 0.6489999890327454
$arr[21,1] = 0.6540567446490791
$arr[21,2] = 0.5186660029865605
$arr[21,3] = 0.6470881035340965
$arr[21,4] = 0.5136884021901443
$arr[22,0] = 0.6330000162124634
$arr[22,1] = 0.6316575410652066
$arr[22,2] = 0.5241413638626182
$arr[22,3] = 0.6326530612244898
$arr[22,4] = 0.5445495271279244
$arr[23,0] = 0.6399999856948853
$arr[23,1] = 0.6346441015430563
$arr[23,2] = 0.5360876057740169
$arr[23,3] = 0.6386261821801892
$arr[23,4] = 0.5286212045793928
$arr[24,0] = 0.6449999809265137
$arr[24,1] = 0.6455948232951717
$arr[24,2] = 0.5415629666500746
$arr[24,3] = 0.6445993031358885
$arr[24,4] = 0.5146839223494276
$arr[25,0] = 0.6365000009536743
$arr[25,1] = 0.6356396217023395
$arr[25,2] = 0.525634644101543
$arr[25,3] = 0.6356396217023395
$arr[25,4] = 0.5296167247386759
$arr[26,0] = 0.6470000147819519
$arr[26,1] = 0.6431060228969636
$arr[26,2] = 0.5181682429069189
$arr[26,3] = 0.6445993031358885
$arr[26,4] = 0.5171727227476356
$arr[27,0] = 0.6420000195503235
$arr[27,1] = 0.6391239422598307
$arr[27,2] = 0.5375808860129417
$arr[27,3] = 0.6396217023394724
$arr[27,4] = 0.5331010452961672
$arr[28,0] = 0.628000020980835
$arr[28,1] = 0.6286709805873569
$arr[28,2] = 0.5201592832254853
$arr[28,3] = 0.6311597809855649
$arr[28,4] = 0.4947735191637631
$arr[29,0] = 0.6395000219345093
$arr[29,1] = 0.6436037829766053
$arr[29,2] = 0.5042309606769537
$arr[29,3] = 0.640119462419114
$arr[29,4] = 0.5126928820308612

$ws.Range("C2:G31").Value = $arr
